$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.296595960855484
$ws.Range("B1").Value = 0.2848637700080872
$ws.Range("C1").Value = 0.2907212972640991
$ws.Range("D1").Value = 0.3799961805343628
$ws.Range("E1").Value = 0.5744473934173584
